$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 467; existing rows 467:498 shift down to 468:499
$ws.Range("A467").EntireRow.Insert()

# Populate the newly inserted row 467 with the new data record
$ws.Range("A467").Value = 6
$ws.Range("B467").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C467").Value = "Metropolitana"
$ws.Range("D467").Value = 44826
$ws.Range("E467").Value = 13
$ws.Range("F467").Value = 100112043
$ws.Range("G467").Value = "Pepino ensalada"
$ws.Range("H467").Value = "Sin especificar"
$ws.Range("I467").Value = "Primera"
$ws.Range("J467").Value = 530
$ws.Range("K467").Value = 16000
$ws.Range("L467").Value = 17000
$ws.Range("M467").Value = 16453
$ws.Range("N467").Value = "`$/caja 60 unidades"
$ws.Range("O467").Value = "Región de Arica y Parinacota"
$ws.Range("P467").Value = 274
$ws.Range("Q467").Value = 60
$ws.Range("R467").Value = "Hortaliza"

# Match the date number format used by the rest of column D
$ws.Range("D467").NumberFormat = "YYYY-MM-DD HH:MM:SS"
